# Fixed update to excel issue
# 1) Rename "Requested quantity" headers on the two existing sheets.
# 2) Add a new "PO Forecast" worksheet (after "Monthly Trend") with the
#    ds / PO_Forecast / yhat_lower / yhat_upper forecast data.

$wb = $excel.ActiveWorkbook

$wsWeekly  = $wb.Worksheets.Item("Weekly Quantity")
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")

$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# Add the new sheet after the last existing sheet ("Monthly Trend").
$lastIndex = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($lastIndex)
$wsForecast = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Mirror the formatting used on "Monthly Trend": bold/centered/bordered
# header row (style index 1) and the date number format on column A
# (style index 2).
$wsMonthly.Range("A1:B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

$wsMonthly.Range("A2").Copy()
$wsForecast.Range("A2:A41").PasteSpecial(-4122)

$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

$data = New-Object 'object[,]' 40,4
    $data[0,0]=45319.99999999999; $data[0,1]=67; $data[0,2]=-766.8595842051694; $data[0,3]=822.675917319228
    $data[1,0]=45326.99999999999; $data[1,1]=98; $data[1,2]=-720.9572726833464; $data[1,3]=943.0050331497337
    $data[2,0]=45333.99999999999; $data[2,1]=130; $data[2,2]=-667.5152406922911; $data[2,3]=931.8595021259689
    $data[3,0]=45340.99999999999; $data[3,1]=161; $data[3,2]=-641.256707932605; $data[3,3]=989.1877945795201
    $data[4,0]=45347.99999999999; $data[4,1]=193; $data[4,2]=-657.4942376451994; $data[4,3]=984.3537540968193
    $data[5,0]=45354.99999999999; $data[5,1]=224; $data[5,2]=-607.4440507959494; $data[5,3]=1006.789410893148
    $data[6,0]=45361.99999999999; $data[6,1]=256; $data[6,2]=-564.729563562996; $data[6,3]=1120.756690539476
    $data[7,0]=45368.99999999999; $data[7,1]=287; $data[7,2]=-575.4790171613984; $data[7,3]=1153.488395318568
    $data[8,0]=45375.99999999999; $data[8,1]=319; $data[8,2]=-482.6010843409741; $data[8,3]=1171.882108543687
    $data[9,0]=45410.99999999999; $data[9,1]=476; $data[9,2]=-367.6528582809731; $data[9,3]=1233.652259103148
    $data[10,0]=45417.99999999999; $data[10,1]=508; $data[10,2]=-277.3058021929779; $data[10,3]=1367.674271769789
    $data[11,0]=45424.99999999999; $data[11,1]=539; $data[11,2]=-250.2071588659605; $data[11,3]=1395.670064100012
    $data[12,0]=45431.99999999999; $data[12,1]=571; $data[12,2]=-180.8314676368104; $data[12,3]=1348.94253979963
    $data[13,0]=45438.99999999999; $data[13,1]=602; $data[13,2]=-155.528716720666; $data[13,3]=1419.597638152831
    $data[14,0]=45445.99999999999; $data[14,1]=634; $data[14,2]=-219.4247787486942; $data[14,3]=1398.233279503762
    $data[15,0]=45459.99999999999; $data[15,1]=697; $data[15,2]=-147.5207648343623; $data[15,3]=1518.317213691385
    $data[16,0]=45466.99999999999; $data[16,1]=729; $data[16,2]=-83.30979620123532; $data[16,3]=1524.438169053004
    $data[17,0]=45473.99999999999; $data[17,1]=760; $data[17,2]=-61.51887162059431; $data[17,3]=1543.327017528362
    $data[18,0]=45487.99999999999; $data[18,1]=823; $data[18,2]=42.73877104514025; $data[18,3]=1618.01102023028
    $data[19,0]=45494.99999999999; $data[19,1]=855; $data[19,2]=5.964126492375379; $data[19,3]=1651.089065287301
    $data[20,0]=45501.99999999999; $data[20,1]=886; $data[20,2]=120.4865946008232; $data[20,3]=1660.647468929117
    $data[21,0]=45508.99999999999; $data[21,1]=918; $data[21,2]=42.09719350108249; $data[21,3]=1740.568072351379
    $data[22,0]=45515.99999999999; $data[22,1]=949; $data[22,2]=70.85831739644409; $data[22,3]=1761.384713313053
    $data[23,0]=45529.99999999999; $data[23,1]=1012; $data[23,2]=196.205485368471; $data[23,3]=1840.508303555779
    $data[24,0]=45536.99999999999; $data[24,1]=1044; $data[24,2]=283.5820787089085; $data[24,3]=1851.32798066132
    $data[25,0]=45543.99999999999; $data[25,1]=1075; $data[25,2]=255.7562636271717; $data[25,3]=1861.775726015381
    $data[26,0]=45550.99999999999; $data[26,1]=1107; $data[26,2]=317.795583148045; $data[26,3]=1946.951685386747
    $data[27,0]=45564.99999999999; $data[27,1]=1170; $data[27,2]=362.8512824070751; $data[27,3]=1961.601633368448
    $data[28,0]=45578.99999999999; $data[28,1]=1233; $data[28,2]=416.4803062176863; $data[28,3]=2037.835777138833
    $data[29,0]=45585.99999999999; $data[29,1]=1265; $data[29,2]=471.691110993772; $data[29,3]=2121.636623119261
    $data[30,0]=45592.99999999999; $data[30,1]=1296; $data[30,2]=447.2025137759381; $data[30,3]=2143.955995903635
    $data[31,0]=45599.99999999999; $data[31,1]=1328; $data[31,2]=552.2753425805778; $data[31,3]=2180.794961413787
    $data[32,0]=45606.99999999999; $data[32,1]=1359; $data[32,2]=560.4650858113642; $data[32,3]=2178.971526127587
    $data[33,0]=45613.99999999999; $data[33,1]=1391; $data[33,2]=557.7575498294461; $data[33,3]=2216.616608121815
    $data[34,0]=45620.99999999999; $data[34,1]=1422; $data[34,2]=698.1711531031513; $data[34,3]=2257.188241981255
    $data[35,0]=45627.99999999999; $data[35,1]=1454; $data[35,2]=571.5315169705906; $data[35,3]=2245.317668582851
    $data[36,0]=45634.99999999999; $data[36,1]=1485; $data[36,2]=661.3805572040477; $data[36,3]=2298.312784880728
    $data[37,0]=45641.99999999999; $data[37,1]=1517; $data[37,2]=656.3815215518237; $data[37,3]=2328.639886526263
    $data[38,0]=45648.99999999999; $data[38,1]=1548; $data[38,2]=687.8828982528911; $data[38,3]=2398.805603487999
    $data[39,0]=45655.99999999999; $data[39,1]=1580; $data[39,2]=804.6858397575421; $data[39,3]=2328.700591956677

$wsForecast.Range("A2:D41").Value = $data

Write-Output "PO Forecast sheet added; headers updated."
